# Update "想去人数" (interest count) figures on the 展览 and 全部类型 sheets.
# The same set of events appears (at different row offsets) on both sheets,
# so locate each row by its current (old) value in column F rather than by
# a fixed row number.
$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Old = 12364; New = 12371 },
    @{ Old = 13;    New = 14 },
    @{ Old = 164;   New = 165 },
    @{ Old = 12194; New = 12200 },
    @{ Old = 4839;  New = 4841 },
    @{ Old = 4718;  New = 4720 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 6)  # Column F
        $val = $cell.Value2
        if ($null -eq $val) { continue }

        foreach ($u in $updates) {
            if ([double]$val -eq [double]$u.Old) {
                $cell.Value = $u.New
                break
            }
        }
    }
}
